# Update the "想去人数" (number of people interested) column F on the
# worksheets "展览", "演出" and "全部类型" to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 522
$ws1.Range("F3").Value = 746
$ws1.Range("F4").Value = 1483
$ws1.Range("F5").Value = 227
$ws1.Range("F6").Value = 92
$ws1.Range("F8").Value = 6198
$ws1.Range("F10").Value = 403
$ws1.Range("F11").Value = 114
$ws1.Range("F12").Value = 5116
$ws1.Range("F14").Value = 178
$ws1.Range("F15").Value = 1176
$ws1.Range("F17").Value = 359
$ws1.Range("F20").Value = 294
$ws1.Range("F22").Value = 3631
$ws1.Range("F23").Value = 151

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 78

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 78
$ws4.Range("F3").Value = 522
$ws4.Range("F4").Value = 746
$ws4.Range("F5").Value = 1483
$ws4.Range("F6").Value = 227
$ws4.Range("F7").Value = 92
$ws4.Range("F9").Value = 6198
$ws4.Range("F11").Value = 403
$ws4.Range("F12").Value = 114
$ws4.Range("F13").Value = 5116
$ws4.Range("F15").Value = 178
$ws4.Range("F16").Value = 1176
$ws4.Range("F18").Value = 359
$ws4.Range("F21").Value = 294
$ws4.Range("F23").Value = 3631
$ws4.Range("F25").Value = 151

$wb.Save()
